# "corrected some formatting errors in flux specific variables"
#
# Row 145 of the "variables-specific" sheet was a stray, half-finished
# "standard_name" metadata row: column B held the shared string
# "standard_name" but column C (the actual value) was left blank. Every
# other variable block in the sheet already has its own, correctly filled
# "standard_name" row immediately below this bogus one, so the blank row
# is pure duplication/clutter. Deleting it shifts all the following rows
# up by one and restores the consistent block layout (dimension shrinks
# from F968 to F967 accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the stray blank row 145 entirely (cells below move up).
$ws.Rows.Item(145).Delete()

# Leave the selection where it was when the fixed file was last saved.
$ws.Range("B80").Select()
